$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the old column I (will no longer be used) ---
$ws.Range("I1:I2").Clear()

# --- Row 1 headers: replace numeric placeholder values with Korean text headers ---
# (A1 previously unused; B1:H1 already carry style s=1 which must be kept)
$ws.Range("A1").Value = "키워드"
$ws.Range("B1").Value = "상품수"
$ws.Range("C1").Value = "한달검색수"
$ws.Range("D1").Value = "6개월매출"
$ws.Range("E1").Value = "6개월판매량"
$ws.Range("F1").Value = "평균가격"
$ws.Range("G1").Value = "경쟁강도"
$ws.Range("H1").Value = "경쟁강도지표"

# A1 had no explicit style before; give it the same header style as the rest of row 1 (s=1)
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows: values that look numeric must still be stored as plain text,
#     except column G which becomes a genuine number. Force "@" (Text) format
#     before writing so Excel does not auto-convert the string to a number,
#     then drop back to the default "Normal" style so no stray style index
#     is left attached to the cell. (Column A holds non-numeric Korean words
#     so it does not need the Text format, avoiding an extra style record.) ---
$ws.Range("B2:F3,H2:H3").NumberFormat = "@"

# A2 previously carried the bordered header style (s=1); the new data row is
# unstyled, so drop it back to the default "Normal" style.
$ws.Range("A2").Style = "Normal"

# Row 2
$ws.Range("A2").Value = "홍당무"
$ws.Range("B2").Value = "493"
$ws.Range("C2").Value = "5567"
$ws.Range("D2").Value = "4500"
$ws.Range("E2").Value = "24310"
$ws.Range("F2").Value = "6190"
$ws.Range("G2").Value = 3.93
$ws.Range("H2").Value = "좋음"

# Row 3 (new row)
$ws.Range("A3").Value = "당근"
$ws.Range("B3").Value = "53129"
$ws.Range("C3").Value = "62581"
$ws.Range("D3").Value = "11800"
$ws.Range("E3").Value = "455787"
$ws.Range("F3").Value = "78300"
$ws.Range("G3").Value = 5.82
$ws.Range("H3").Value = "좋음"

# Reset to the default (no explicit) style now that the text is safely stored
$ws.Range("B2:F3,H2:H3").Style = "Normal"
